$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.222.82'
$ws.Range('E2').Value = '  -1.61%  '

$ws.Range('D3').Value = '2.959.59'
$ws.Range('E3').Value = '  -3.36%  '

$ws.Range('E4').Value = '  -0.23%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '497.35'
$ws.Range('E5').Value = '  -2.80%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.79'
$ws.Range('E6').Value = '  +4.19%  '

$ws.Range('E7').Value = '  -0.22%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.425'
$ws.Range('E8').Value = '  -1.18%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.23'
$ws.Range('E9').Value = '  +2.41%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.105'
$ws.Range('E10').Value = '  +1.41%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.348'
$ws.Range('E11').Value = '  -3.05%  '

$ws.Range('E12').Value = '  +0.08%  '

$ws.Range('D13').Value = '3.465.48'
$ws.Range('E13').Value = '  -4.30%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.08'
$ws.Range('E14').Value = '  +3.74%  '

$ws.Range('D15').Value = '56.212.90'
$ws.Range('E15').Value = '  +2.82%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000147'
$ws.Range('E16').Value = '  +2.52%  '

$ws.Range('D17').Value = '2.960.10'
$ws.Range('E17').Value = '  -4.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.66'
$ws.Range('E18').Value = '  +2.91%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.25'
$ws.Range('E19').Value = '  -1.33%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.69'
$ws.Range('E20').Value = '  +1.76%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.50'
$ws.Range('E21').Value = '  -1.21%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.05%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.467'
$ws.Range('E23').Value = '  -4.32%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.72'
$ws.Range('E24').Value = '  -5.70%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.11%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.161'
$ws.Range('E26').Value = '  -2.58%  '

$ws.Range('D27').Value = '0.0₃0877'
$ws.Range('E27').Value = '  +0.59%  '

$ws.Range('E28').Value = '  -0.11%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.38'
$ws.Range('E29').Value = '  -1.81%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.73'
$ws.Range('E30').Value = '  +2.30%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.72'
$ws.Range('E31').Value = '  -3.55%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  -5.22%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.11'
$ws.Range('E33').Value = '  -2.41%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '154.73'
$ws.Range('E34').Value = '  -0.90%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.40'
$ws.Range('E35').Value = '  -3.48%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.27'
$ws.Range('E36').Value = '  -2.70%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.50'
$ws.Range('E37').Value = '  -6.33%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0665'
$ws.Range('E38').Value = '  +1.71%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.60'
$ws.Range('E39').Value = '  -1.11%  '

$ws.Range('D40').Value = '2.990.69'
$ws.Range('E40').Value = '  -3.91%  '

$ws.Range('E41').Value = '  -0.12%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '35.77'
$ws.Range('E42').Value = '  -1.17%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.637'
$ws.Range('E43').Value = '  -3.92%  '

$ws.Range('D44').Value = '2.217.02'
$ws.Range('E44').Value = '  -0.08%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.979'
$ws.Range('E45').Value = '  -5.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.37'
$ws.Range('E46').Value = '  +0.44%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.52'
$ws.Range('E47').Value = '  -4.86%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0233'
$ws.Range('E48').Value = '  +3.26%  '

$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.89'
$ws.Range('E49').Value = '  +10.40%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.71'
$ws.Range('E50').Value = '  -3.66%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.70'
$ws.Range('E51').Value = '  -4.12%  '
